$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 9000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 9000
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -9338

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 414.6
$ws.Range("I33").Value = 433.66666
$ws.Range("J33").Value = 357.4
$ws.Range("K33").Value = 433.66666
$ws.Range("L33").Value = 357.4
$ws.Range("M33").Value = -204.66666
$ws.Range("N33").Value = -815.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 425.46155
$ws.Range("I38").Value = 169.25
$ws.Range("K38").Value = 507.75
$ws.Range("M38").Value = -135.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 175.91667
$ws.Range("I39").Value = 146
$ws.Range("K39").Value = 438
$ws.Range("M39").Value = -142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9300.4
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4825

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 292.18182
$ws.Range("I42").Value = 228.55556
$ws.Range("J42").Value = 578.5
$ws.Range("K42").Value = 685.66668
$ws.Range("L42").Value = 1735.5
$ws.Range("M42").Value = -455.66668
$ws.Range("N42").Value = -2195.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 2523.182
$ws.Range("I45").Value = 856.7143
$ws.Range("K45").Value = 2570.1429
$ws.Range("M45").Value = -2378.1429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8667.666999999999
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 12001.5
$ws.Range("K64").Value = 2000
$ws.Range("L64").Value = 12001.5
$ws.Range("M64").Value = -1752
$ws.Range("N64").Value = -12497.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 8667.666999999999
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 12001.5
$ws.Range("K67").Value = 2000
$ws.Range("L67").Value = 12001.5
$ws.Range("M67").Value = -1142
$ws.Range("N67").Value = -13717.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 13482
$ws.Range("I74").Value = 9970.6
$ws.Range("J74").Value = 19334.334
$ws.Range("K74").Value = 9970.6
$ws.Range("L74").Value = 19334.334
$ws.Range("M74").Value = -9034.6
$ws.Range("N74").Value = -21206.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6496.8184
$ws.Range("I76").Value = 5247
$ws.Range("K76").Value = 5247
$ws.Range("M76").Value = -4932

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 13482
$ws.Range("I77").Value = 9970.6
$ws.Range("J77").Value = 19334.334
$ws.Range("K77").Value = 49853
$ws.Range("L77").Value = 96671.67
$ws.Range("M77").Value = -45173
$ws.Range("N77").Value = -106031.67

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6496.8184
$ws.Range("I79").Value = 5247
$ws.Range("K79").Value = 5247
$ws.Range("M79").Value = -4155

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1499.5
$ws.Range("I111").Value = 1499.5
$ws.Range("K111").Value = 4498.5
$ws.Range("M111").Value = -1431.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1669.7307
$ws.Range("J112").Value = 1783.4546
$ws.Range("L112").Value = 5350.3638
$ws.Range("N112").Value = -7566.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3897.052
$ws.Range("I32").Value = 2945.6807
$ws.Range("K32").Value = 2945.6807
$ws.Range("M32").Value = -2658.6807

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3199.6
$ws.Range("I122").Value = 2557.5
$ws.Range("K122").Value = 7672.5
$ws.Range("M122").Value = -5222.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3600.1428
$ws.Range("I132").Value = 2695.4783
$ws.Range("J132").Value = 5334.0835
$ws.Range("K132").Value = 8086.4349
$ws.Range("L132").Value = 16002.2505
$ws.Range("M132").Value = -5556.4349
$ws.Range("N132").Value = -21062.2505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 69990.336
$ws.Range("J21").Value = 69990.336
$ws.Range("L21").Value = 69990.336
$ws.Range("N21").Value = -70462.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 291
$ws.Range("I22").Value = 291
$ws.Range("K22").Value = 291
$ws.Range("M22").Value = -118

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 25000
$ws.Range("J46").Value = 25000
$ws.Range("L46").Value = 25000
$ws.Range("N46").Value = -25596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 20666.666
$ws.Range("J49").Value = 20666.666
$ws.Range("L49").Value = 20666.666
$ws.Range("N49").Value = -21144.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3240.3333
$ws.Range("I86").Value = 3126.2666
$ws.Range("K86").Value = 3126.2666
$ws.Range("M86").Value = -2003.2666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3240.3333
$ws.Range("I89").Value = 3126.2666
$ws.Range("K89").Value = 15631.333
$ws.Range("M89").Value = -10015.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1331.8667
$ws.Range("I107").Value = 855.5714
$ws.Range("K107").Value = 855.5714
$ws.Range("M107").Value = 1064.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3790.3333
$ws.Range("I134").Value = 2157
$ws.Range("K134").Value = 6471
$ws.Range("M134").Value = -3936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1223.28
$ws.Range("I16").Value = 629.05884
$ws.Range("J16").Value = 2486
$ws.Range("K16").Value = 629.05884
$ws.Range("L16").Value = 2486
$ws.Range("M16").Value = -342.05884
$ws.Range("N16").Value = -3060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6494.8
$ws.Range("I22").Value = 4589.8
$ws.Range("K22").Value = 4589.8
$ws.Range("M22").Value = -4239.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26465.637
$ws.Range("I31").Value = 2022.0869
$ws.Range("J31").Value = 53237.145
$ws.Range("K31").Value = 2022.0869
$ws.Range("L31").Value = 53237.145
$ws.Range("M31").Value = -1727.0869
$ws.Range("N31").Value = -53827.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 26465.637
$ws.Range("I34").Value = 2022.0869
$ws.Range("J34").Value = 53237.145
$ws.Range("K34").Value = 2022.0869
$ws.Range("L34").Value = 53237.145
$ws.Range("M34").Value = -1820.0869
$ws.Range("N34").Value = -53641.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3691.5625
$ws.Range("I58").Value = 1821.4615
$ws.Range("J58").Value = 11795.333
$ws.Range("K58").Value = 1821.4615
$ws.Range("L58").Value = 11795.333
$ws.Range("M58").Value = -1618.4615
$ws.Range("N58").Value = -12201.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9219.700000000001
$ws.Range("I62").Value = 4031.8333
$ws.Range("K62").Value = 4031.8333
$ws.Range("M62").Value = -3407.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 9219.700000000001
$ws.Range("I65").Value = 4031.8333
$ws.Range("K65").Value = 20159.1665
$ws.Range("M65").Value = -17039.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 15664.833
$ws.Range("I93").Value = 15664.833
$ws.Range("K93").Value = 15664.833
$ws.Range("M93").Value = -13792.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1223.28
$ws.Range("I113").Value = 629.05884
$ws.Range("J113").Value = 2486
$ws.Range("K113").Value = 629.05884
$ws.Range("L113").Value = 2486
$ws.Range("M113").Value = 1540.94116
$ws.Range("N113").Value = -6826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3937.5789
$ws.Range("I134").Value = 2688.3635
$ws.Range("J134").Value = 5655.25
$ws.Range("K134").Value = 8065.0905
$ws.Range("L134").Value = 16965.75
$ws.Range("M134").Value = -5530.0905
$ws.Range("N134").Value = -22035.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 69677.39999999999
$ws.Range("J135").Value = 69677.39999999999
$ws.Range("L135").Value = 69677.39999999999
$ws.Range("N135").Value = -79817.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3691.5625
$ws.Range("I136").Value = 1821.4615
$ws.Range("J136").Value = 11795.333
$ws.Range("K136").Value = 5464.3845
$ws.Range("L136").Value = 35385.999
$ws.Range("M136").Value = -2914.3845
$ws.Range("N136").Value = -40485.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 5000
$ws.Range("K14").Value = 15000
$ws.Range("M14").Value = -14827

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 8500.666999999999
$ws.Range("J80").Value = 9000.5
$ws.Range("L80").Value = 27001.5
$ws.Range("N80").Value = -28873.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 8500.666999999999
$ws.Range("J83").Value = 9000.5
$ws.Range("L83").Value = 81004.5
$ws.Range("N83").Value = -90364.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1751.5
$ws.Range("I97").Value = 1751.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 5254.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -4758.5
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4449
$ws.Range("J104").Value = 4449
$ws.Range("L104").Value = 13347
$ws.Range("N104").Value = -18589

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1286.1666
$ws.Range("I121").Value = 1316.875
$ws.Range("J121").Value = 1224.75
$ws.Range("K121").Value = 3950.625
$ws.Range("L121").Value = 3674.25
$ws.Range("M121").Value = -2640.625
$ws.Range("N121").Value = -6294.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 41670770
$ws.Range("J124").Value = 55560210
$ws.Range("L124").Value = 166680630
$ws.Range("N124").Value = -166690450

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 843.3077
$ws.Range("I2").Value = 161.33333
$ws.Range("J2").Value = 1427.8572
$ws.Range("K2").Value = 161.33333
$ws.Range("L2").Value = 1427.8572
$ws.Range("M2").Value = -48.33332999999999
$ws.Range("N2").Value = -1653.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 50000
$ws.Range("J15").Value = 50000
$ws.Range("L15").Value = 50000
$ws.Range("N15").Value = -50576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 39999
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 14280.75
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 39999
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 50000
$ws.Range("J57").Value = 50000
$ws.Range("L57").Value = 50000
$ws.Range("N57").Value = -51640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 459163
$ws.Range("I80").Value = 558420.75
$ws.Range("K80").Value = 558420.75
$ws.Range("M80").Value = -557422.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 459163
$ws.Range("I83").Value = 558420.75
$ws.Range("K83").Value = 2792103.75
$ws.Range("M83").Value = -2787111.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1655.1351
$ws.Range("I102").Value = 984.2593000000001
$ws.Range("K102").Value = 984.2593000000001
$ws.Range("M102").Value = 637.7406999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2915.8235
$ws.Range("I132").Value = 2360.2917
$ws.Range("K132").Value = 7080.875100000001
$ws.Range("M132").Value = -4550.875100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 70525.78999999999
$ws.Range("J133").Value = 70525.78999999999
$ws.Range("L133").Value = 70525.78999999999
$ws.Range("N133").Value = -80645.78999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3425.1853
$ws.Range("J22").Value = 4681.091
$ws.Range("L22").Value = 4681.091
$ws.Range("N22").Value = -5271.091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3425.1853
$ws.Range("J27").Value = 4681.091
$ws.Range("L27").Value = 4681.091
$ws.Range("N27").Value = -4895.091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8999.762000000001
$ws.Range("I46").Value = 6748.25
$ws.Range("K46").Value = 6748.25
$ws.Range("M46").Value = -6560.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1668972.9
$ws.Range("I55").Value = 3572952.8
$ws.Range("J55").Value = 2990.5
$ws.Range("K55").Value = 3572952.8
$ws.Range("L55").Value = 2990.5
$ws.Range("M55").Value = -3572779.8
$ws.Range("N55").Value = -3336.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -101232

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2084.7693
$ws.Range("I93").Value = 2093.5833
$ws.Range("K93").Value = 2093.5833
$ws.Range("M93").Value = -845.5832999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 84890
$ws.Range("J130").Value = 84890
$ws.Range("L130").Value = 84890
$ws.Range("N130").Value = -94930

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5852.6
$ws.Range("I132").Value = 6072.647
$ws.Range("J132").Value = 5172.4546
$ws.Range("K132").Value = 18217.941
$ws.Range("L132").Value = 15517.3638
$ws.Range("M132").Value = -15687.941
$ws.Range("N132").Value = -20577.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 42662.8
$ws.Range("I70").Value = 40000
$ws.Range("K70").Value = 40000
$ws.Range("M70").Value = -39685

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 42662.8
$ws.Range("I73").Value = 40000
$ws.Range("K73").Value = 40000
$ws.Range("M73").Value = -38908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 49998.332
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 49998.332
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2093.5117
$ws.Range("I122").Value = 1454.2424
$ws.Range("K122").Value = 4362.7272
$ws.Range("M122").Value = -1912.7272

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 71955.664
$ws.Range("J130").Value = 71955.664
$ws.Range("L130").Value = 71955.664
$ws.Range("N130").Value = -81995.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3232.9795
$ws.Range("I136").Value = 2057.4167
$ws.Range("K136").Value = 6172.250100000001
$ws.Range("M136").Value = -3622.250100000001
